$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(-16.95087480019198,2.078630820641369,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,2.160487598133653,2.264831188974931,-16.95087480019198,3.44916343847627,-16.95087480019198,2.034725355757201,-16.95087480019198,2.356897454784786,-16.95087480019198),
    @(-16.95087480019198,1.897494458074535,-16.95087480019198,-16.95087480019198,-16.95087480019198,3.008641863039727,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(2.59007240920921,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,1.737378924408372,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(3.805144744361071,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,2.048610535927522),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,2.858078380066515,-16.95087480019198,2.418349777859334,-16.95087480019198,-16.95087480019198,-16.95087480019198,1.344037566746744),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,2.418775790828863,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,2.190087242177887,1.878665233148658),
    @(-16.95087480019198,-16.95087480019198,1.183335909894623,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,2.134683239703399),
    @(-16.95087480019198,-16.95087480019198,1.264928338321209,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,2.347127337650412,-16.95087480019198),
    @(-16.95087480019198,1.941979256099772,2.266820261025026,-16.95087480019198,-16.95087480019198,-16.95087480019198,1.210360438109745,-16.95087480019198,1.402758804495896,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,-16.95087480019198,1.10915083348889,-16.95087480019198,1.368325411467133,-16.95087480019198),
    @(-16.95087480019198,-16.95087480019198,1.502753536936703,-16.95087480019198,-16.95087480019198,-16.95087480019198,1.467772566458414,4.321917276114449,-16.95087480019198,-16.95087480019198),
    @(-16.95087480019198,0.7919411851529271,1.522708088111276,-16.95087480019198,3.182353974927274,-16.95087480019198,1.918613617784502,-16.95087480019198,-16.95087480019198,2.393813373805465),
    @(-16.95087480019198,1.004167727291364,-16.95087480019198,2.024293050635714,-16.95087480019198,2.723822645579043,2.288287830439973,-16.95087480019198,-16.95087480019198,-16.95087480019198)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 2 + $i
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
